$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price/volume snapshot (and the Dai/ShibaInu
# row-order swap) as scraped on Thu Jun  8 14:08:03 UTC 2023.
#
# Every Price/Volume(1h) cell in the sheet is stored as plain TEXT, even
# when the text happens to look like a number (e.g. "1.001"). Writing such
# a value straight into Range.Value would make Excel auto-convert it to a
# real number, so those cells are written with a leading apostrophe to force
# text, then the quote-prefix style that leaves behind is cleared via
# Style = "Normal" so the cell keeps the same (default) style as before.

$ws.Range("D2").Value = "26.384.57"
$ws.Range("E2").Value = "  -1.54%  "
$ws.Range("D3").Value = "1.845.93"
$ws.Range("E3").Value = "  -0.92%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'264.34"
$ws.Range("E5").Value = "  -1.48%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "'0.5196"
$ws.Range("E7").Value = "  -1.84%  "
$ws.Range("D8").Value = "'0.3268"
$ws.Range("E8").Value = "  -2.36%  "
$ws.Range("D9").Value = "'0.06803"
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("D10").Value = "'18.76"
$ws.Range("E10").Value = "  -4.54%  "
$ws.Range("D11").Value = "'0.7753"
$ws.Range("E11").Value = "  -0.75%  "
$ws.Range("D12").Value = "'0.07774"
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("D13").Value = "1.844.89"
$ws.Range("E13").Value = "  -1.18%  "
$ws.Range("D14").Value = "'87.93"
$ws.Range("E14").Value = "  -2.24%  "
$ws.Range("D15").Value = "'5.011"
$ws.Range("E15").Value = "  -1.89%  "
$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").Value = "'13.93"
$ws.Range("E17").Value = "  -2.92%  "
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").Value = "'1.001"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.000007968"
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("D20").Value = "26.418.79"
$ws.Range("E20").Value = "  -1.51%  "
$ws.Range("D21").Value = "2.075.92"
$ws.Range("E21").Value = "  -1.76%  "
$ws.Range("D22").Value = "'4.643"
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").Value = "'9.537"
$ws.Range("E23").Value = "  -3.21%  "
$ws.Range("D24").Value = "'5.983"
$ws.Range("E24").Value = "  -1.06%  "
$ws.Range("D25").Value = "'144.67"
$ws.Range("E25").Value = "  -0.57%  "
$ws.Range("D26").Value = "'2.210"
$ws.Range("E26").Value = "  -8.24%  "
$ws.Range("D27").Value = "'1.660"
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("D28").Value = "'17.02"
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("D29").Value = "'112.18"
$ws.Range("E29").Value = "  -0.50%  "
$ws.Range("D30").Value = "'4.180"
$ws.Range("E30").Value = "  -2.32%  "
$ws.Range("D31").Value = "'4.134"
$ws.Range("E31").Value = "  -3.05%  "
$ws.Range("D32").Value = "'0.08750"
$ws.Range("D33").Value = "'0.04830"
$ws.Range("E33").Value = "  -2.15%  "
$ws.Range("D34").Value = "'1.135"
$ws.Range("E34").Value = "  -1.68%  "
$ws.Range("D35").Value = "'0.7188"
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("D36").Value = "'2.860"
$ws.Range("E36").Value = "  -0.74%  "
$ws.Range("D37").Value = "'3.097"
$ws.Range("E37").Value = "  -2.66%  "
$ws.Range("D38").Value = "'0.01778"
$ws.Range("D39").Value = "'2.191"
$ws.Range("E39").Value = "  -4.37%  "
$ws.Range("D40").Value = "'0.4855"
$ws.Range("E40").Value = "  -3.29%  "
$ws.Range("D41").Value = "'0.9295"
$ws.Range("E41").Value = "  +1.51%  "
$ws.Range("D42").Value = "'111.00"
$ws.Range("E42").Value = "  -3.92%  "
$ws.Range("D43").Value = "'6.057"
$ws.Range("E43").Value = "  -1.00%  "
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "'7.693"
$ws.Range("E45").Value = "  -2.93%  "
$ws.Range("D46").Value = "'0.05933"
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("D47").Value = "'0.4160"
$ws.Range("E47").Value = "  -4.47%  "
$ws.Range("D48").Value = "'9.110"
$ws.Range("E48").Value = "  -1.97%  "
$ws.Range("D49").Value = "'0.1236"
$ws.Range("E49").Value = "  -6.00%  "
$ws.Range("D50").Value = "'34.90"
$ws.Range("E50").Value = "  -2.55%  "
$ws.Range("D51").Value = "'0.8911"
$ws.Range("E51").Value = "  +1.57%  "

# Clear the forced-text (quote-prefix) style hint cell by cell so every
# touched cell keeps style 0, same as the rest of the data rows.
$textForcedRefs = @(
    "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D17", "D18", "D19", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D50", "D51"
)
foreach ($ref in $textForcedRefs) {
    $ws.Range($ref).Style = "Normal"
}
